$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cases tab query (B2): the "Cohort" output column/line was dropped ---
$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.stage_of_disease IN [''Va'']
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $casesQuery

# --- StatQuery column (C2:C4): replaced with the new filtered/aliquot-aware stats query ---
$statsQuery = ' MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])
        AND (s.study_disposition = ''Unrestricted'')
        AND (size([]) = 0 OR s.clinical_study_type IN [])
        AND (size([''Va'']) = 0 OR demo.breed IN [''Va''])
        AND (size([]) = 0 OR demo.sex IN [])
        AND (size([]) = 0 OR demo.neutered_indicator IN [])
        AND (size([]) = 0 OR diag.disease_term IN [])
        AND (size([]) = 0 OR diag.primary_disease_site IN [])
        AND (size([]) = 0 OR diag.stage_of_disease IN [])
        AND (size([]) = 0 OR diag.best_response IN [])
    OPTIONAL MATCH (c)-->(co:cohort)
    OPTIONAL MATCH (f:file)-[*]->(c)
    OPTIONAL MATCH (f)-->(parent)
    OPTIONAL MATCH (samp:sample)-->(c)
    OPTIONAL MATCH (samp)<--(al:aliquot)
    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al
      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])
        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])
        AND (size([]) = 0 OR samp.sample_site IN [])
        AND (size([]) = 0 OR head(labels(parent)) IN [])
        AND (size([]) = 0 OR f.file_type IN [])
        AND (size([]) = 0 OR f.file_format IN [])
    WITH c.case_id AS case_id,
         s.clinical_study_designation AS study_code,
         s.clinical_study_type AS study_type,
         co.cohort_description AS cohort,
         demo.breed AS breed,
         diag.disease_term AS diagnosis,
         diag.stage_of_disease AS stage_of_disease,
         diag.primary_disease_site AS disease_site,
         demo.patient_age_at_enrollment AS age,
         demo.sex AS sex,
         demo.neutered_indicator AS neutered_status,
         demo.weight AS weight,
         diag.best_response AS response_to_treatment,
         samp.sample_id AS sample_id,
         f.uuid AS file_id,
         al
    RETURN
COUNT(DISTINCT file_id) as number_of_files,
COUNT(DISTINCT sample_id) as number_of_sample,
COUNT(DISTINCT case_id) as number_of_cases,
COUNT(DISTINCT study_code) as number_of_study,
COUNT(DISTINCT al) as number_of_aliquot'
$ws.Range("C2").Value = $statsQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("C4").Value = $statsQuery

# --- Cosmetic re-save side effects (zoom level bump + rows growing to the
#     max wrap height once the StatQuery text got a lot longer) ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6
